# Weekly update: insert a new daily record for the carrot ("Zanahoria")
# price series. A new row is inserted immediately above the current
# row 167, pushing all the subsequent rows (old 167..251) down by one
# (they become 168..252), and the newly inserted row 167 is populated
# with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 167 - this shifts rows 167:251
# down to 168:252 (and grows the used range from R251 to R252).
$ws.Rows("167:167").Insert()

# Populate the newly inserted row 167 with the new weekly record.
$ws.Cells.Item(167, 1).Value  = 1
$ws.Cells.Item(167, 2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(167, 3).Value  = 'Arica y Parinacota'
$ws.Cells.Item(167, 4).Value  = 44636
$ws.Cells.Item(167, 5).Value  = 15
$ws.Cells.Item(167, 6).Value  = 100114013
$ws.Cells.Item(167, 7).Value  = 'Zanahoria'
$ws.Cells.Item(167, 8).Value  = 'Sin especificar'
$ws.Cells.Item(167, 9).Value  = 'Primera'
$ws.Cells.Item(167, 10).Value = 70
$ws.Cells.Item(167, 11).Value = 18000
$ws.Cells.Item(167, 12).Value = 19000
$ws.Cells.Item(167, 13).Value = 18500
$ws.Cells.Item(167, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(167, 15).Value = 'Valle de Camiña'
$ws.Cells.Item(167, 16).Value = 740
$ws.Cells.Item(167, 17).Value = 25
$ws.Cells.Item(167, 18).Value = 'Hortaliza'
